$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table data: headers (row 1) + 3 measure rows (rows 2-4), columns A-E
# ---------------------------------------------------------------------------
$rows = @(
    @("Measure ID", "Measure Name", "Measure Status", "Measure Description", "Applicable Specialties"),
    @("IA_PM_25", "Standardized Risk Assessment", "New", "Evidence-based risk stratification", "Cardiology, Primary Care"),
    @("IA_PM_26", "Lifestyle Intervention Support", "Updated", "Updated lifestyle coaching documentation", "Internal Medicine, Endocrinology"),
    @("IA_EPA_1", "Obsolete Measure Example", "Deleted", "Deprecated due to overlap with IA_PM_25", "All Specialties")
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowValues = $rows[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# ---------------------------------------------------------------------------
# Header row (row 1) formatting: bold font, thin box border, centered
# horizontally and top-aligned vertically.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Restore the cell selection to match the saved workbook state.
# ---------------------------------------------------------------------------
$null = $ws.Range("D14").Select()
